$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.072.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.719.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.45%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.712.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.645"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.774"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000408"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +58.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.300.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.19%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.704.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.60%  "
$ws.Range("E20").Value = "  +4.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.152.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +24.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("E31").Value = "  +9.38%  "
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0748"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +33.11%  "
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "29.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +35.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("E45").Value = "  +33.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("E47").Value = "  +5.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.310"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.23%  "
